$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the average formula to exclude C2 (history score) from the range,
# now averaging only C4:C9.
$ws.Range("C10").Formula = "=AVERAGE(C4:C9)"

# Update the active selection on Sheet1 to F11.
$ws.Range("F11").Select()
